# Applies the diff:
#  1. Removes the "Meta description: ..." paragraph near the top of the
#     document (it used to follow the "Play Danger High Voltage Megapays
#     Free!" heading).
#  2. Inserts a new paragraph - "Play Danger High Voltage Megapays Free!"
#     in bold - right before the final (italic) paragraph.
#  3. Replaces the text of that final paragraph (the old AI image prompt)
#     with the meta-description sentence that used to live near the top.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph ----------------------
# Find the paragraph whose text starts with "Meta description" and remove
# the whole paragraph (including its own paragraph mark) so the document
# collapses back together cleanly.
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete() | Out-Null
}

# --- Step 2: insert a new bold paragraph before the last paragraph --------
# Insert a brand-new, clean paragraph after the second-to-last paragraph
# (avoids inheriting the trailing paragraph's italic formatting / list
# style), then fill it in with the bold heading text.
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$r = $secondLast.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter([char]13)

$newPara = $d.Paragraphs.Item($count)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Danger High Voltage Megapays Free!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xml) | Out-Null

# --- Step 3: swap the old AI image prompt for the meta description text ---
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for Danger High Voltage Megapays that features a happy Maya warrior wearing a pair of stylish glasses. The warrior should be standing in front of a disco ball with bright lights shining behind them. The image should be vibrant and eye-catching, with the warrior looking confident and ready to take on the reels. Add some electric effects around the edges of the image to highlight the " + [char]34 + "Danger" + [char]34 + " aspect of the game's title. Overall, the image should be fun, lively, and capture the essence of the game's upbeat soundtrack and dancefloor theme.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Danger High Voltage Megapays, play for free, and learn how to win progressive jackpots!",
    2) | Out-Null
